$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1878980891719745
$ws.Range("C2").Value = 0.5636942675159236
$ws.Range("J2").Value = 0.03184713375796178
$ws.Range("P2").Value = 0.1146496815286624
$ws.Range("S2").Value = 0.1019108280254777
$ws.Range("B3").Value = 0.01052631578947368
$ws.Range("C3").Value = 0.06315789473684211
$ws.Range("J3").Value = 0.06842105263157895
$ws.Range("P3").Value = 0.7
$ws.Range("S3").Value = 0.1578947368421053
$ws.Range("J4").Value = 0.075
$ws.Range("P4").Value = 0.575
$ws.Range("S4").Value = 0.35
$ws.Range("B6").Value = 0.06779661016949153
$ws.Range("D6").Value = 0.02542372881355932
$ws.Range("E6").Value = 0.00423728813559322
$ws.Range("F6").Value = 0.07203389830508475
$ws.Range("J6").Value = 0.2711864406779661
$ws.Range("O6").Value = 0.01271186440677966
$ws.Range("Q6").Value = 0.1610169491525424
$ws.Range("R6").Value = 0.05932203389830509
$ws.Range("S6").Value = 0.326271186440678
$ws.Range("B7").Value = 0.09547738693467336
$ws.Range("D7").Value = 0.005025125628140704
$ws.Range("F7").Value = 0.04020100502512563
$ws.Range("J7").Value = 0.08040201005025126
$ws.Range("O7").Value = 0.03015075376884422
$ws.Range("Q7").Value = 0.2412060301507538
$ws.Range("R7").Value = 0.1005025125628141
$ws.Range("S7").Value = 0.407035175879397
$ws.Range("B8").Value = 0.1120162932790224
$ws.Range("D8").Value = 0.01425661914460285
$ws.Range("F8").Value = 0.07331975560081466
$ws.Range("J8").Value = 0.07942973523421588
$ws.Range("O8").Value = 0.01629327902240326
$ws.Range("Q8").Value = 0.2016293279022403
$ws.Range("R8").Value = 0.09979633401221996
$ws.Range("S8").Value = 0.4032586558044807
$ws.Range("B9").Value = 0.1256830601092896
$ws.Range("D9").Value = 0.0273224043715847
$ws.Range("F9").Value = 0.06557377049180328
$ws.Range("J9").Value = 0.07650273224043716
$ws.Range("O9").Value = 0.01092896174863388
$ws.Range("Q9").Value = 0.1530054644808743
$ws.Range("R9").Value = 0.1147540983606557
$ws.Range("S9").Value = 0.4262295081967213
$ws.Range("B10").Value = 0.09819494584837545
$ws.Range("D10").Value = 0.01732851985559567
$ws.Range("E10").Value = 0.001444043321299639
$ws.Range("F10").Value = 0.05342960288808665
$ws.Range("J10").Value = 0.09458483754512635
$ws.Range("O10").Value = 0.01227436823104693
$ws.Range("Q10").Value = 0.2223826714801444
$ws.Range("R10").Value = 0.1155234657039711
$ws.Range("S10").Value = 0.3848375451263538
$ws.Range("F11").Value = 0.003095975232198143
$ws.Range("G11").Value = 0.1331269349845201
$ws.Range("J11").Value = 0.09907120743034056
$ws.Range("K11").Value = 0.2074303405572755
$ws.Range("L11").Value = 0.5325077399380805
$ws.Range("S11").Value = 0.02476780185758514
$ws.Range("G12").Value = 0.7329545454545454
$ws.Range("J12").Value = 0.1988636363636364
$ws.Range("K12").Value = 0.01704545454545454
$ws.Range("L12").Value = 0.03409090909090909
$ws.Range("S12").Value = 0.01704545454545454
$ws.Range("G13").Value = 0.6862745098039216
$ws.Range("J13").Value = 0.2941176470588235
$ws.Range("S13").Value = 0.0196078431372549
$ws.Range("J14").Value = 0.5
$ws.Range("S14").Value = 0.5
$ws.Range("F15").Value = 0.03125
$ws.Range("H15").Value = 0.171875
$ws.Range("I15").Value = 0.046875
$ws.Range("J15").Value = 0.37890625
$ws.Range("K15").Value = 0.046875
$ws.Range("M15").Value = 0.0078125
$ws.Range("N15").Value = 0.00390625
$ws.Range("O15").Value = 0.0859375
$ws.Range("S15").Value = 0.2265625
$ws.Range("F16").Value = 0.03243243243243243
$ws.Range("H16").Value = 0.1675675675675676
$ws.Range("I16").Value = 0.07567567567567568
$ws.Range("J16").Value = 0.3567567567567568
$ws.Range("K16").Value = 0.1189189189189189
$ws.Range("M16").Value = 0.02702702702702703
$ws.Range("O16").Value = 0.1027027027027027
$ws.Range("S16").Value = 0.1189189189189189
$ws.Range("F17").Value = 0.03281853281853282
$ws.Range("H17").Value = 0.1563706563706564
$ws.Range("I17").Value = 0.07915057915057915
$ws.Range("J17").Value = 0.4517374517374517
$ws.Range("K17").Value = 0.08880308880308881
$ws.Range("M17").Value = 0.02509652509652509
$ws.Range("O17").Value = 0.06756756756756757
$ws.Range("S17").Value = 0.09845559845559845
$ws.Range("F18").Value = 0.04562737642585551
$ws.Range("H18").Value = 0.1673003802281369
$ws.Range("I18").Value = 0.06844106463878327
$ws.Range("J18").Value = 0.4296577946768061
$ws.Range("K18").Value = 0.09885931558935361
$ws.Range("M18").Value = 0.02281368821292776
$ws.Range("O18").Value = 0.07224334600760456
$ws.Range("S18").Value = 0.09505703422053231
$ws.Range("F19").Value = 0.01186063750926612
$ws.Range("H19").Value = 0.2157153446997776
$ws.Range("I19").Value = 0.07412898443291327
$ws.Range("J19").Value = 0.3839881393624907
$ws.Range("K19").Value = 0.1052631578947368
$ws.Range("M19").Value = 0.02075611564121571
$ws.Range("N19").Value = 0.001482579688658265
$ws.Range("O19").Value = 0.07116382505559674
$ws.Range("S19").Value = 0.1156412157153447
